{"js": "// Fix contact information missing from short resumes:\n// Insert a new, centered paragraph with the contact info directly\n// after the \"Dheeraj Chand\" title paragraph (the first paragraph in\n// the document body), matching the target OOXML exactly (a plain\n// <w:r><w:t>...</w:t></w:r> run with no direct run formatting, inside\n// a paragraph whose <w:pPr> only carries <w:jc w:val=\"center\"/>).\n//\n// A plain Paragraph.insertParagraph()/insertText() call would inherit\n// the bold/28-half-point direct run formatting of the \"Dheeraj Chand\"\n// run it is anchored next to, which does not match the target XML (no\n// <w:rPr> at all). Inserting a literal OOXML fragment sidesteps that\n// formatting-inheritance behavior and reproduces the target markup\n// precisely.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0];\nconst titleRange = titleParagraph.getRange();\n\nconst contactText =\n  \"202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX\";\n\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n  '<w:r><w:t>' + contactText + '</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntitleRange.insertOoxml(ooxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Fix contact information missing from short resumes:\n# Insert a new, centered paragraph with the contact info directly\n# after the \"Dheeraj Chand\" title paragraph (the document's first\n# paragraph), matching the target OOXML exactly: a plain\n# <w:r><w:t>...</w:t></w:r> run (no direct run formatting) inside a\n# paragraph whose <w:pPr> only carries <w:jc w:val=\"center\"/>.\n#\n# A straightforward Range.InsertParagraphAfter() + Range.Text=... would\n# inherit the bold / 28-half-point direct run formatting of the\n# \"Dheeraj Chand\" run it is anchored to, which would not match the\n# target XML (no <w:rPr> at all on the new run). Using Find/Replace with\n# a literal paragraph mark (^p) in the replacement text instead performs\n# a true paragraph split at that point, and the newly created run picks\n# up no direct character formatting - exactly matching the target\n# markup. The paragraph-level alignment (<w:jc w:val=\"center\"/>) carries\n# forward automatically from the \"Dheeraj Chand\" paragraph's own\n# properties, which is also what the target shows.\n\n$d = $word.ActiveDocument\n\n$contactText = \"202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Dheeraj Chand\"\n$find.Replacement.Text = \"Dheeraj Chand^p\" + $contactText\n\n$find.Execute(\n    $find.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find.Replacement.Text,\n    2\n)\n"}
